# Scheduled runner update: refresh currentAveragePrice-derived columns
# (H, I, J, K, L, M, N) across the Leve profit sheets with new market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 20837084
$ws.Range("J51").Value = 20837084
$ws.Range("L51").Value = 20837084
$ws.Range("N51").Value = -20838052
$ws.Range("H57").Value = 55842.332
$ws.Range("J57").Value = 55842.332
$ws.Range("L57").Value = 167526.996
$ws.Range("N57").Value = -168524.996
$ws.Range("H74").Value = 3246.8948
$ws.Range("I74").Value = 2477.9285
$ws.Range("K74").Value = 2477.9285
$ws.Range("M74").Value = -1541.9285
$ws.Range("H76").Value = 3686.4
$ws.Range("I76").Value = 3695.75
$ws.Range("J76").Value = 3649
$ws.Range("K76").Value = 3695.75
$ws.Range("L76").Value = 3649
$ws.Range("M76").Value = -3380.75
$ws.Range("N76").Value = -4279
$ws.Range("H77").Value = 3246.8948
$ws.Range("I77").Value = 2477.9285
$ws.Range("K77").Value = 12389.6425
$ws.Range("M77").Value = -7709.6425
$ws.Range("H79").Value = 3686.4
$ws.Range("I79").Value = 3695.75
$ws.Range("J79").Value = 3649
$ws.Range("K79").Value = 3695.75
$ws.Range("L79").Value = 3649
$ws.Range("M79").Value = -2603.75
$ws.Range("N79").Value = -5833
$ws.Range("H137").Value = 1282.25
$ws.Range("I137").Value = 1226
$ws.Range("K137").Value = 3678
$ws.Range("M137").Value = -1128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 20377.25
$ws.Range("I28").Value = 16996.666
$ws.Range("J28").Value = 30519
$ws.Range("K28").Value = 16996.666
$ws.Range("L28").Value = 30519
$ws.Range("M28").Value = -16804.666
$ws.Range("N28").Value = -30903
$ws.Range("H61").Value = 1804.1666
$ws.Range("I61").Value = 1443.5
$ws.Range("K61").Value = 1443.5
$ws.Range("M61").Value = -1231.5
$ws.Range("H74").Value = 1869
$ws.Range("I74").Value = 1869
$ws.Range("K74").Value = 1869
$ws.Range("M74").Value = -995
$ws.Range("H77").Value = 1869
$ws.Range("I77").Value = 1869
$ws.Range("K77").Value = 9345
$ws.Range("M77").Value = -4977
$ws.Range("H99").Value = 20377.25
$ws.Range("I99").Value = 16996.666
$ws.Range("J99").Value = 30519
$ws.Range("K99").Value = 16996.666
$ws.Range("L99").Value = 30519
$ws.Range("M99").Value = -14001.666
$ws.Range("N99").Value = -36509
$ws.Range("H110").Value = 1096.2941
$ws.Range("I110").Value = 1150.1333
$ws.Range("J110").Value = 692.5
$ws.Range("K110").Value = 1150.1333
$ws.Range("L110").Value = 692.5
$ws.Range("M110").Value = 894.8667
$ws.Range("N110").Value = -4782.5
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H136").Value = 1804.1666
$ws.Range("I136").Value = 1443.5
$ws.Range("K136").Value = 4330.5
$ws.Range("M136").Value = -1780.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 89997
$ws.Range("J2").Value = 89997
$ws.Range("L2").Value = 89997
$ws.Range("N2").Value = -90223
$ws.Range("H6").Value = 22601
$ws.Range("J6").Value = 22601
$ws.Range("L6").Value = 22601
$ws.Range("N6").Value = -22827
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H134").Value = 1558.7391
$ws.Range("I134").Value = 1422.6
$ws.Range("J134").Value = 2466.3333
$ws.Range("K134").Value = 4267.799999999999
$ws.Range("L134").Value = 7398.999899999999
$ws.Range("M134").Value = -1732.799999999999
$ws.Range("N134").Value = -12468.9999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 95459.37
$ws.Range("J122").Value = 9026.200000000001
$ws.Range("L122").Value = 27078.6
$ws.Range("N122").Value = -31978.6
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("K132").Value = 18000
$ws.Range("M132").Value = -15470
$ws.Range("H134").Value = 4628.615
$ws.Range("I134").Value = 2799.8333
$ws.Range("J134").Value = 6196.143
$ws.Range("K134").Value = 8399.499899999999
$ws.Range("L134").Value = 18588.429
$ws.Range("M134").Value = -5864.499899999999
$ws.Range("N134").Value = -23658.429
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 706.625
$ws.Range("I11").Value = 195.8
$ws.Range("K11").Value = 587.4000000000001
$ws.Range("M11").Value = -447.4000000000001
$ws.Range("H64").Value = 2813.8125
$ws.Range("J64").Value = 3256.125
$ws.Range("L64").Value = 9768.375
$ws.Range("N64").Value = -10308.375
$ws.Range("H67").Value = 2813.8125
$ws.Range("J67").Value = 3256.125
$ws.Range("L67").Value = 9768.375
$ws.Range("N67").Value = -11640.375
$ws.Range("H98").Value = 712.44446
$ws.Range("I98").Value = 862
$ws.Range("J98").Value = 413.33334
$ws.Range("K98").Value = 2586
$ws.Range("L98").Value = 1240.00002
$ws.Range("M98").Value = -1088
$ws.Range("N98").Value = -4236.000019999999
$ws.Range("H116").Value = 4995
$ws.Range("I116").Value = 4995
$ws.Range("K116").Value = 14985
$ws.Range("M116").Value = -11543
$ws.Range("H131").Value = 1667.6666
$ws.Range("J131").Value = 1689.75
$ws.Range("L131").Value = 5069.25
$ws.Range("N131").Value = -15149.25
$ws.Range("H132").Value = 1785.9584
$ws.Range("I132").Value = 963.8889
$ws.Range("J132").Value = 2279.2
$ws.Range("K132").Value = 8675.000100000001
$ws.Range("L132").Value = 20512.8
$ws.Range("M132").Value = -6145.000100000001
$ws.Range("N132").Value = -25572.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 12900
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H113").Value = 3742.3809
$ws.Range("I113").Value = 3408.75
$ws.Range("K113").Value = 3408.75
$ws.Range("M113").Value = -1238.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 18148.814
$ws.Range("I46").Value = 34309.23
$ws.Range("K46").Value = 34309.23
$ws.Range("M46").Value = -34121.23
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19132
$ws.Range("H61").Value = 80882.94
$ws.Range("I61").Value = 73334.14
$ws.Range("K61").Value = 73334.14
$ws.Range("M61").Value = -73132.14
$ws.Range("H113").Value = 80882.94
$ws.Range("I113").Value = 73334.14
$ws.Range("K113").Value = 73334.14
$ws.Range("M113").Value = -71164.14
$ws.Range("H131").Value = 84987
$ws.Range("J131").Value = 84987
$ws.Range("L131").Value = 84987
$ws.Range("N131").Value = -95067
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 17080
$ws.Range("J14").Value = 40500
$ws.Range("L14").Value = 40500
$ws.Range("N14").Value = -40836
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H126").Value = 2016.5834
$ws.Range("I126").Value = 1982.1818
$ws.Range("K126").Value = 5946.5454
$ws.Range("M126").Value = -3476.5454
$ws.Range("H130").Value = 49333.332
$ws.Range("J130").Value = 49333.332
$ws.Range("L130").Value = 49333.332
$ws.Range("N130").Value = -59373.332
